# Insert two new rows at the top of the data block (row 985), pushing the
# existing historical rows down by two, and fill the new rows with the new
# week's price data for "Femacal de La Calera" / Pina / Caramelo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 985..1101 down by two rows.
$ws.Rows.Item(985).Resize(2).Insert()

# New row 985: "Primera" quality, 108 units/caja, $23000/caja 12 unidades
$ws.Range("A985").Value = 3
$ws.Range("B985").Value = "Femacal de La Calera"
$ws.Range("C985").Value = "Coquimbo"
$ws.Range("D985").Value = 45142
$ws.Range("E985").Value = 5
$ws.Range("F985").Value = "Fruta"
$ws.Range("G985").Value = 100108
$ws.Range("H985").Value = "Tropicales y subtropicales"
$ws.Range("I985").Value = 100108005
$ws.Range("J985").Value = "Piña"
$ws.Range("K985").Value = "Caramelo"
$ws.Range("L985").Value = "Primera"
$ws.Range("M985").Value = 108
$ws.Range("N985").Value = 23000
$ws.Range("O985").Value = 23000
$ws.Range("P985").Value = 23000
$ws.Range("Q985").Value = "$/caja 12 unidades"
$ws.Range("R985").Value = "Ecuador"
$ws.Range("S985").Value = 1917
$ws.Range("T985").Value = 12

# New row 986: "Segunda" quality, 108 units/caja, $23000/caja 14 unidades
$ws.Range("A986").Value = 3
$ws.Range("B986").Value = "Femacal de La Calera"
$ws.Range("C986").Value = "Coquimbo"
$ws.Range("D986").Value = 45142
$ws.Range("E986").Value = 5
$ws.Range("F986").Value = "Fruta"
$ws.Range("G986").Value = 100108
$ws.Range("H986").Value = "Tropicales y subtropicales"
$ws.Range("I986").Value = 100108005
$ws.Range("J986").Value = "Piña"
$ws.Range("K986").Value = "Caramelo"
$ws.Range("L986").Value = "Segunda"
$ws.Range("M986").Value = 108
$ws.Range("N986").Value = 23000
$ws.Range("O986").Value = 23000
$ws.Range("P986").Value = 23000
$ws.Range("Q986").Value = "$/caja 14 unidades"
$ws.Range("R986").Value = "Ecuador"
$ws.Range("S986").Value = 1643
$ws.Range("T986").Value = 14
